$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 64 (Leve Item ID G=5506)
$ws_ALC.Cells.Item(64, 8).Value = 4082.7
$ws_ALC.Cells.Item(64, 9).Value = 3555
$ws_ALC.Cells.Item(64, 10).Value = 4214.625
$ws_ALC.Cells.Item(64, 11).Value = 3555
$ws_ALC.Cells.Item(64, 12).Value = 4214.625
$ws_ALC.Cells.Item(64, 13).Value = -3307
$ws_ALC.Cells.Item(64, 14).Value = -4710.625

# ALC row 67 (Leve Item ID G=5506)
$ws_ALC.Cells.Item(67, 8).Value = 4082.7
$ws_ALC.Cells.Item(67, 9).Value = 3555
$ws_ALC.Cells.Item(67, 10).Value = 4214.625
$ws_ALC.Cells.Item(67, 11).Value = 3555
$ws_ALC.Cells.Item(67, 12).Value = 4214.625
$ws_ALC.Cells.Item(67, 13).Value = -2697
$ws_ALC.Cells.Item(67, 14).Value = -5930.625

# ALC row 70 (Leve Item ID G=12604)
$ws_ALC.Cells.Item(70, 8).Value = 3003
$ws_ALC.Cells.Item(70, 10).Value = 3003
$ws_ALC.Cells.Item(70, 12).Value = 9009
$ws_ALC.Cells.Item(70, 14).Value = -9549

# ALC row 73 (Leve Item ID G=12604)
$ws_ALC.Cells.Item(73, 8).Value = 3003
$ws_ALC.Cells.Item(73, 10).Value = 3003
$ws_ALC.Cells.Item(73, 12).Value = 9009
$ws_ALC.Cells.Item(73, 14).Value = -10881

# ALC row 74 (Leve Item ID G=5507)
$ws_ALC.Cells.Item(74, 8).Value = 4236.55
$ws_ALC.Cells.Item(74, 9).Value = 4523
$ws_ALC.Cells.Item(74, 11).Value = 4523
$ws_ALC.Cells.Item(74, 13).Value = -3587

# ALC row 77 (Leve Item ID G=5507)
$ws_ALC.Cells.Item(77, 8).Value = 4236.55
$ws_ALC.Cells.Item(77, 9).Value = 4523
$ws_ALC.Cells.Item(77, 11).Value = 22615
$ws_ALC.Cells.Item(77, 13).Value = -17935

# ARM row 61 (Leve Item ID G=43999)
$ws_ARM.Cells.Item(61, 8).Value = 91093336
$ws_ARM.Cells.Item(61, 9).Value = 125127336
$ws_ARM.Cells.Item(61, 10).Value = 336001.34
$ws_ARM.Cells.Item(61, 11).Value = 125127336
$ws_ARM.Cells.Item(61, 12).Value = 336001.34
$ws_ARM.Cells.Item(61, 13).Value = -125127124
$ws_ARM.Cells.Item(61, 14).Value = -336425.34

# ARM row 124 (Leve Item ID G=34252)
$ws_ARM.Cells.Item(124, 8).Value = 21439.5
$ws_ARM.Cells.Item(124, 10).Value = 21439.5
$ws_ARM.Cells.Item(124, 12).Value = 21439.5
$ws_ARM.Cells.Item(124, 14).Value = -31259.5

# ARM row 132 (Leve Item ID G=43997)
$ws_ARM.Cells.Item(132, 8).Value = 46122.223
$ws_ARM.Cells.Item(132, 9).Value = 32755.438
$ws_ARM.Cells.Item(132, 10).Value = 79025.08
$ws_ARM.Cells.Item(132, 11).Value = 98266.314
$ws_ARM.Cells.Item(132, 12).Value = 237075.24
$ws_ARM.Cells.Item(132, 13).Value = -95736.314
$ws_ARM.Cells.Item(132, 14).Value = -242135.24

# ARM row 136 (Leve Item ID G=43999)
$ws_ARM.Cells.Item(136, 8).Value = 91093336
$ws_ARM.Cells.Item(136, 9).Value = 125127336
$ws_ARM.Cells.Item(136, 10).Value = 336001.34
$ws_ARM.Cells.Item(136, 11).Value = 375382008
$ws_ARM.Cells.Item(136, 12).Value = 1008004.02
$ws_ARM.Cells.Item(136, 13).Value = -375379458
$ws_ARM.Cells.Item(136, 14).Value = -1013104.02

# BSM row 22 (Leve Item ID G=5092)
$ws_BSM.Cells.Item(22, 8).Value = 732.7692
$ws_BSM.Cells.Item(22, 9).Value = 567.125
$ws_BSM.Cells.Item(22, 11).Value = 567.125
$ws_BSM.Cells.Item(22, 13).Value = -394.125

# CRP row 58 (Leve Item ID G=44021)
$ws_CRP.Cells.Item(58, 8).Value = 41670240
$ws_CRP.Cells.Item(58, 9).Value = 41670240
$ws_CRP.Cells.Item(58, 10).Value = 0
$ws_CRP.Cells.Item(58, 11).Value = 41670240
$ws_CRP.Cells.Item(58, 12).Value = 0
$ws_CRP.Cells.Item(58, 13).Value = $null
$ws_CRP.Cells.Item(58, 14).Value = -41670037

# CRP row 62 (Leve Item ID G=12580)
$ws_CRP.Cells.Item(62, 8).Value = 0
$ws_CRP.Cells.Item(62, 9).Value = 0
$ws_CRP.Cells.Item(62, 10).Value = 0
$ws_CRP.Cells.Item(62, 11).Value = 0
$ws_CRP.Cells.Item(62, 12).Value = $null
$ws_CRP.Cells.Item(62, 13).Value = $null
$ws_CRP.Cells.Item(62, 14).Value = 0

# CRP row 65 (Leve Item ID G=12580)
$ws_CRP.Cells.Item(65, 8).Value = 0
$ws_CRP.Cells.Item(65, 9).Value = 0
$ws_CRP.Cells.Item(65, 10).Value = 0
$ws_CRP.Cells.Item(65, 11).Value = 0
$ws_CRP.Cells.Item(65, 12).Value = $null
$ws_CRP.Cells.Item(65, 13).Value = $null
$ws_CRP.Cells.Item(65, 14).Value = 0

# CRP row 99 (Leve Item ID G=36198)
$ws_CRP.Cells.Item(99, 8).Value = 4617.5293
$ws_CRP.Cells.Item(99, 9).Value = 4617.5293
$ws_CRP.Cells.Item(99, 11).Value = 4617.5293
$ws_CRP.Cells.Item(99, 13).Value = -3119.5293

# CRP row 122 (Leve Item ID G=36196)
$ws_CRP.Cells.Item(122, 8).Value = 1453.8379
$ws_CRP.Cells.Item(122, 9).Value = 1174.2142
$ws_CRP.Cells.Item(122, 10).Value = 2323.7778
$ws_CRP.Cells.Item(122, 11).Value = 3522.6426
$ws_CRP.Cells.Item(122, 12).Value = 6971.3334
$ws_CRP.Cells.Item(122, 13).Value = -1072.6426
$ws_CRP.Cells.Item(122, 14).Value = -11871.3334

# CRP row 126 (Leve Item ID G=36198)
$ws_CRP.Cells.Item(126, 8).Value = 4617.5293
$ws_CRP.Cells.Item(126, 9).Value = 4617.5293
$ws_CRP.Cells.Item(126, 11).Value = 13852.5879
$ws_CRP.Cells.Item(126, 13).Value = -11382.5879

# CRP row 136 (Leve Item ID G=44021)
$ws_CRP.Cells.Item(136, 8).Value = 41670240
$ws_CRP.Cells.Item(136, 9).Value = 41670240
$ws_CRP.Cells.Item(136, 10).Value = 0
$ws_CRP.Cells.Item(136, 11).Value = 125010720
$ws_CRP.Cells.Item(136, 12).Value = 0
$ws_CRP.Cells.Item(136, 13).Value = -125008170
$ws_CRP.Cells.Item(136, 14).Value = -125008170

# CUL row 92 (Leve Item ID G=19841)
$ws_CUL.Cells.Item(92, 8).Value = 787.8
$ws_CUL.Cells.Item(92, 10).Value = 997
$ws_CUL.Cells.Item(92, 12).Value = 2991
$ws_CUL.Cells.Item(92, 14).Value = -5487

# CUL row 131 (Leve Item ID G=36060)
$ws_CUL.Cells.Item(131, 8).Value = 924.1951
$ws_CUL.Cells.Item(131, 10).Value = 978.1622
$ws_CUL.Cells.Item(131, 12).Value = 2934.4866
$ws_CUL.Cells.Item(131, 14).Value = -13014.4866

# GSM row 70 (Leve Item ID G=14146)
$ws_GSM.Cells.Item(70, 8).Value = 5065.276
$ws_GSM.Cells.Item(70, 9).Value = 4725.9473
$ws_GSM.Cells.Item(70, 11).Value = 4725.9473
$ws_GSM.Cells.Item(70, 13).Value = -4455.9473

# GSM row 73 (Leve Item ID G=14146)
$ws_GSM.Cells.Item(73, 8).Value = 5065.276
$ws_GSM.Cells.Item(73, 9).Value = 4725.9473
$ws_GSM.Cells.Item(73, 11).Value = 4725.9473
$ws_GSM.Cells.Item(73, 13).Value = -3789.9473

# GSM row 102 (Leve Item ID G=36169)
$ws_GSM.Cells.Item(102, 8).Value = 1995.1538
$ws_GSM.Cells.Item(102, 9).Value = 1402.8
$ws_GSM.Cells.Item(102, 11).Value = 1402.8
$ws_GSM.Cells.Item(102, 13).Value = 219.2

# GSM row 122 (Leve Item ID G=36182)
$ws_GSM.Cells.Item(122, 8).Value = 2709.5833
$ws_GSM.Cells.Item(122, 9).Value = 2166.5557
$ws_GSM.Cells.Item(122, 10).Value = 4338.6665
$ws_GSM.Cells.Item(122, 11).Value = 6499.6671
$ws_GSM.Cells.Item(122, 12).Value = 13015.9995
$ws_GSM.Cells.Item(122, 13).Value = -4049.6671
$ws_GSM.Cells.Item(122, 14).Value = -17915.9995

# GSM row 132 (Leve Item ID G=44008)
$ws_GSM.Cells.Item(132, 8).Value = 100058.81
$ws_GSM.Cells.Item(132, 9).Value = 88231.914
$ws_GSM.Cells.Item(132, 10).Value = 115828
$ws_GSM.Cells.Item(132, 11).Value = 264695.742
$ws_GSM.Cells.Item(132, 12).Value = 347484
$ws_GSM.Cells.Item(132, 13).Value = -262165.742
$ws_GSM.Cells.Item(132, 14).Value = -352544

# LTW row 40 (Leve Item ID G=36248)
$ws_LTW.Cells.Item(40, 8).Value = 7142.857
$ws_LTW.Cells.Item(40, 9).Value = 7142.857
$ws_LTW.Cells.Item(40, 11).Value = 7142.857
$ws_LTW.Cells.Item(40, 13).Value = -7006.857

# LTW row 55 (Leve Item ID G=5284)
$ws_LTW.Cells.Item(55, 8).Value = 223.44897
$ws_LTW.Cells.Item(55, 9).Value = 197.45714
$ws_LTW.Cells.Item(55, 10).Value = 288.42856
$ws_LTW.Cells.Item(55, 11).Value = 197.45714
$ws_LTW.Cells.Item(55, 12).Value = 288.42856
$ws_LTW.Cells.Item(55, 13).Value = -24.45714000000001
$ws_LTW.Cells.Item(55, 14).Value = -634.4285600000001

# LTW row 76 (Leve Item ID G=10742)
$ws_LTW.Cells.Item(76, 8).Value = 36513.816
$ws_LTW.Cells.Item(76, 10).Value = 36513.816
$ws_LTW.Cells.Item(76, 12).Value = 36513.816
$ws_LTW.Cells.Item(76, 14).Value = -37189.816

# LTW row 79 (Leve Item ID G=10742)
$ws_LTW.Cells.Item(79, 8).Value = 36513.816
$ws_LTW.Cells.Item(79, 10).Value = 36513.816
$ws_LTW.Cells.Item(79, 12).Value = 36513.816
$ws_LTW.Cells.Item(79, 14).Value = -38853.816

# LTW row 132 (Leve Item ID G=44058)
$ws_LTW.Cells.Item(132, 8).Value = 87649.836
$ws_LTW.Cells.Item(132, 9).Value = 3850
$ws_LTW.Cells.Item(132, 10).Value = 171449.67
$ws_LTW.Cells.Item(132, 11).Value = 11550
$ws_LTW.Cells.Item(132, 12).Value = 514349.01
$ws_LTW.Cells.Item(132, 13).Value = -9020
$ws_LTW.Cells.Item(132, 14).Value = -519409.01

# LTW row 136 (Leve Item ID G=44060)
$ws_LTW.Cells.Item(136, 8).Value = 168866.25
$ws_LTW.Cells.Item(136, 9).Value = 126875.625
$ws_LTW.Cells.Item(136, 10).Value = 252847.5
$ws_LTW.Cells.Item(136, 11).Value = 380626.875
$ws_LTW.Cells.Item(136, 12).Value = 758542.5
$ws_LTW.Cells.Item(136, 13).Value = -378076.875
$ws_LTW.Cells.Item(136, 14).Value = -763642.5

# WVR row 100 (Leve Item ID G=19981)
$ws_WVR.Cells.Item(100, 8).Value = 56280.5
$ws_WVR.Cells.Item(100, 9).Value = 42262.418
$ws_WVR.Cells.Item(100, 10).Value = 84316.664
$ws_WVR.Cells.Item(100, 11).Value = 84524.836
$ws_WVR.Cells.Item(100, 12).Value = 168633.328
$ws_WVR.Cells.Item(100, 13).Value = -83983.836
$ws_WVR.Cells.Item(100, 14).Value = -169715.328

# WVR row 132 (Leve Item ID G=44029)
$ws_WVR.Cells.Item(132, 8).Value = 119303.06
$ws_WVR.Cells.Item(132, 9).Value = 111939.336
$ws_WVR.Cells.Item(132, 10).Value = 127587.25
$ws_WVR.Cells.Item(132, 11).Value = 335818.008
$ws_WVR.Cells.Item(132, 12).Value = 382761.75
$ws_WVR.Cells.Item(132, 13).Value = -333288.008
$ws_WVR.Cells.Item(132, 14).Value = -387821.75

# WVR row 136 (Leve Item ID G=44031)
$ws_WVR.Cells.Item(136, 8).Value = 121382.3
$ws_WVR.Cells.Item(136, 9).Value = 87482.836
$ws_WVR.Cells.Item(136, 10).Value = 202741
$ws_WVR.Cells.Item(136, 11).Value = 262448.508
$ws_WVR.Cells.Item(136, 12).Value = 608223
$ws_WVR.Cells.Item(136, 13).Value = -259898.508
$ws_WVR.Cells.Item(136, 14).Value = -613323
